$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the W/L/T values for each data row
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 96
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 0
}
